$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values for the duplicate_image_filename column (E) for rows 2 through 21
$ws.Range("E2:E21").Value = "NA"
